$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.07"
$ws.Range("E2").Value = "'4.79%"
$ws.Range("D3").Value = "'36.01"
$ws.Range("E3").Value = "'15.90%"
$ws.Range("E4").Value = "'4.53%"
$ws.Range("D5").Value = "'0.07869"
$ws.Range("E5").Value = "'6.94%"
$ws.Range("D6").Value = "'2.324"
$ws.Range("E6").Value = "'1.17%"
$ws.Range("D7").Value = "'8.056"
$ws.Range("D8").Value = "'3.972"
$ws.Range("E8").Value = "'6.33%"
$ws.Range("D9").Value = "'0.9265"
$ws.Range("E9").Value = "'0.90%"
$ws.Range("D10").Value = "'0.1012"
$ws.Range("E10").Value = "'10.36%"
$ws.Range("D11").Value = "'0.1821"
$ws.Range("E11").Value = "'7.11%"
$ws.Range("D12").Value = "'0.08481"
$ws.Range("E12").Value = "'2.10%"
$ws.Range("D13").Value = "'0.03373"
$ws.Range("E13").Value = "'8.38%"
$ws.Range("D14").Value = "'0.09914"
$ws.Range("E14").Value = "'-0.77%"
$ws.Range("D15").Value = "'0.001467"
$ws.Range("E15").Value = "'-1.90%"
$ws.Range("D16").Value = "'0.005777"
$ws.Range("E16").Value = "'0.51%"
$ws.Range("D17").Value = "'3.481"
$ws.Range("E17").Value = "'0.26%"
$ws.Range("D18").Value = "'2.127"
$ws.Range("E18").Value = "'1.13%"
$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'3.04%"
$ws.Range("D21").Value = "'4.530"
$ws.Range("E21").Value = "'8.55%"
$ws.Range("D22").Value = "'0.2216"
$ws.Range("E22").Value = "'4.37%"
$ws.Range("D23").Value = "'0.04628"
$ws.Range("E23").Value = "'2.91%"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'0.04%"
$ws.Range("D25").Value = "'0.004464"
$ws.Range("E25").Value = "'6.39%"
$ws.Range("E26").Value = "'-0.30%"
$ws.Range("D27").Value = "'0.0003389"
$ws.Range("E27").Value = "'-0.12%"
$ws.Range("D39").Value = "'0.01758"
$ws.Range("E39").Value = "'11.71%"
$ws.Range("D40").Value = "'0.04738"
$ws.Range("E40").Value = "'5.27%"
$ws.Range("D41").Value = "'0.007910"
$ws.Range("E41").Value = "'7.38%"
$ws.Range("D42").Value = "'0.1416"
$ws.Range("E42").Value = "'5.87%"
$ws.Range("D43").Value = "'0.008797"
$ws.Range("E43").Value = "'-10.65%"
$ws.Range("D44").Value = "'0.002211"
$ws.Range("E44").Value = "'-0.38%"
$ws.Range("D45").Value = "'0.009153"
$ws.Range("E45").Value = "'7.48%"
$ws.Range("D46").Value = "'0.00006063"
$ws.Range("E46").Value = "'-0.66%"
$ws.Range("D47").Value = "'0.00000000747"
$ws.Range("E47").Value = "'-0.31%"
$ws.Range("D48").Value = "'3.904"
$ws.Range("E48").Value = "'59.64%"
$ws.Range("D49").Value = "'0.002682"
$ws.Range("E49").Value = "'34.18%"
$ws.Range("D50").Value = "'0.00002093"
$ws.Range("E50").Value = "'-0.31%"
$ws.Range("D51").Value = "'0.0001993"
$ws.Range("E51").Value = "'-0.31%"
